$d = $word.ActiveDocument

# Locate the "ОПЫТ РАБОТЫ" heading run inside the info table.
$rng = $d.Content
$found = $rng.Find.Execute("ОПЫТ РАБОТЫ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find 'ОПЫТ РАБОТЫ' heading text"
}

# Collapse to the end of the match, then append " (2+ лет)" right after it.
$rng.Collapse(0)
$rng.InsertAfter(" (2+ лет)")

# The inserted text currently shares the exact same run properties as the
# heading it follows, so the engine merges it into a single run. Nudge the
# color away and then restore it to the original value (262626) so the new
# text keeps its own separate <w:r> run (matching formatting, distinct run)
# instead of being silently absorbed into the preceding run.
$rng.Font.Color = 255
$rng.Font.Color = 2500134
